# Update the document per the commit diff:
# - date line
# - 25 three-digit division problems

$d = $word.ActiveDocument

$replacements = @(
    @{old = "2026-01-19 Monday"; new = "2026-01-20 Tuesday"},
    @{old = "501÷9="; new = "903÷5="},
    @{old = "597÷2="; new = "977÷5="},
    @{old = "745÷8="; new = "684÷8="},
    @{old = "340÷8="; new = "497÷2="},
    @{old = "318÷5="; new = "659÷6="},
    @{old = "814÷6="; new = "265÷4="},
    @{old = "848÷7="; new = "471÷6="},
    @{old = "113÷8="; new = "150÷3="},
    @{old = "995÷3="; new = "821÷2="},
    @{old = "583÷4="; new = "787÷4="},
    @{old = "880÷3="; new = "412÷2="},
    @{old = "288÷4="; new = "466÷3="},
    @{old = "275÷5="; new = "589÷9="},
    @{old = "397÷5="; new = "785÷9="},
    @{old = "620÷4="; new = "631÷3="},
    @{old = "802÷7="; new = "742÷6="},
    @{old = "582÷5="; new = "994÷5="},
    @{old = "431÷9="; new = "357÷2="},
    @{old = "782÷6="; new = "191÷6="},
    @{old = "761÷8="; new = "759÷8="},
    @{old = "667÷5="; new = "282÷7="},
    @{old = "229÷6="; new = "280÷2="},
    @{old = "910÷5="; new = "453÷6="},
    @{old = "666÷4="; new = "218÷9="},
    @{old = "152÷4="; new = "492÷9="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
